$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text (matching the original inlineStr formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "26.083.96"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "1.598.60"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  +2.81%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("D9").Value = "0.0616"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "17.96"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  +4.43%  "
$ws.Range("D12").Value = "1.819.91"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "1.594.54"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "26.060.54"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "60.41"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "204.70"
$ws.Range("E20").Value = "  +10.78%  "
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").Value = "9.31"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "5.98"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  +11.81%  "
$ws.Range("D25").Value = "141.22"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").Value = "15.22"
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("D29").Value = "6.45"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").Value = "1.109.25"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").Value = "0.0162"
$ws.Range("E37").Value = "  +8.72%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "0.776"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").Value = "0.494"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").Value = "1.732.21"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.10"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "92.62"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +5.17%  "
$ws.Range("D47").Value = "53.35"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "0.0502"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "0.408"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "7.23"
$ws.Range("E51").Value = "  +1.04%  "
